$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A so the existing "reference",
# "document date" and "gross amount" columns (A,B,C) shift right to C,D,E.
$ws.Range("A1:B1").EntireColumn.Insert()

# New header row values for the two inserted columns.
$ws.Cells.Item(1, 1).Value = "vendor"
$ws.Cells.Item(1, 2).Value = "doc. number"

# New vendor / doc. number data for rows 2-7.
$vendor = @(5461563, 654632, 654631, 6546323, 6865, 987651)
$docNumber = @(78000163, 78000185, 78000207, 78000229, 78000251, 78000273)

for ($i = 0; $i -lt $vendor.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $vendor[$i]
    $ws.Cells.Item($row, 2).Value = $docNumber[$i]
}

# Width for the new "doc. number" column (column B) — closest value this
# host's ColumnWidth setter (quantized to 1/6 character units) can reach to
# the authored 14.5703125.
$ws.Columns.Item(2).ColumnWidth = 13.67

# Update the selection to match the new layout.
$ws.Range("F6").Select() | Out-Null
